$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -------------------------------------------------------------
# Row 2 ("stop" button) count: 2 -> 10
$ws.Range("B2").Value = 10

# Row 4 ("reset" button): barcode text + count change
$ws.Range("A4").Value = "a978895940627"
$ws.Range("B4").Value = 5

# --- Style edits --------------------------------------------------------
# The barcode-label column (A2:A4) switches from the themed/minor font with
# an explicit left+center alignment to a plain "Malgun Gothic" font with
# just the inherited vertical-center alignment (no horizontal override).
$labels = $ws.Range("A2:A4")
$labels.ClearFormats()
$labels.Font.Name = "Malgun Gothic"

# --- View state ---------------------------------------------------------
# Move the active selection.
$ws.Range("B13").Select()
